# Add new data row 38 to the "Artfynd" worksheet, mirroring the existing
# row layout/style (no explicit cell styles are used anywhere in this sheet).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 38

function Set-TextValue($r, $c, $text) {
    # Force a literal text value even when it looks like a number/date/time
    # (e.g. "1", "2023-06-26", "00:00") so Excel doesn't silently coerce it
    # into a number/date serial. Resetting the style back to "Normal"
    # afterwards drops the transient quote-prefix / text-format styling so
    # the cell ends up using the shared default style, matching the rest
    # of the sheet (which has no per-cell styles at all).
    $cell = $ws.Cells.Item($r, $c)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

# A38 Id
$ws.Cells.Item($row, 1).Value = 112182731
# B38 Taxonsorteringsordning
$ws.Cells.Item($row, 2).Value = 89965
# C38 Valideringsstatus
$ws.Cells.Item($row, 3).Value = "Ovaliderad"
# D38 Rodlistade
$ws.Cells.Item($row, 4).Value = "VU"
# E38 TaxonId
$ws.Cells.Item($row, 5).Value = 760
# F38 Artnamn
$ws.Cells.Item($row, 6).Value = "Doftticka"
# G38 Vetenskapligt namn
$ws.Cells.Item($row, 7).Value = "Haploporus odorus"
# H38 Auktor
$ws.Cells.Item($row, 8).Value = "(Sommerf.) Bondartsev & Singer"
# I38 Antal (text "1", not numeric 1)
Set-TextValue $row 9 "1"
# P38 Lokalnamn
$ws.Cells.Item($row, 16).Value = "Lycksaberg, Ly lm"
# Q38 Ost
$ws.Cells.Item($row, 17).Value = 636505.8934863589
# R38 Nord
$ws.Cells.Item($row, 18).Value = 7215899.030802801
# S38 Noggrannhet
$ws.Cells.Item($row, 19).Value = 5
# T38 Lan
$ws.Cells.Item($row, 20).Value = "Västerbotten"
# U38 Kommun
$ws.Cells.Item($row, 21).Value = "Storuman"
# V38 Provins
$ws.Cells.Item($row, 22).Value = "Lycksele lappmark"
# W38 Forsamling
$ws.Cells.Item($row, 23).Value = "Stensele"
# Y38 Startdatum (text, not a date serial)
Set-TextValue $row 25 "2023-06-26"
# Z38 Starttid (text, not a time serial)
Set-TextValue $row 26 "00:00"
# AA38 Slutdatum (text, not a date serial)
Set-TextValue $row 27 "2023-06-26"
# AB38 Sluttid (text, not a time serial)
Set-TextValue $row 28 "00:00"
# AC38 Publik kommentar
$ws.Cells.Item($row, 29).Value = "Påträffad under Sveaskogs naturvärdesinventering"
# AD38 Ej återfunnen
$ws.Cells.Item($row, 30).Value = $false
# AE38 Osäker artbestämning
$ws.Cells.Item($row, 31).Value = $false
# AG38 Ospontan
$ws.Cells.Item($row, 33).Value = $false
# AT38 Bestämningsår (blank placeholder column, present but empty)
$ws.Cells.Item($row, 46).Formula = '=""'
# AW38 Rapportör
$ws.Cells.Item($row, 49).Value = "Mimmi Persson"
# AX38 Observatörer
$ws.Cells.Item($row, 50).Value = "Mimmi Persson"
# AY38 Projektnamn (blank placeholder column, present but empty)
$ws.Cells.Item($row, 51).Formula = '=""'
